$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 124, shifting existing rows 124:244 down to 125:245
$ws.Rows(124).Insert()

# Populate the new row 124 with the new data record
$ws.Range("A124").Value = 4
$ws.Range("B124").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C124").Value = "Los Lagos"
$ws.Range("D124").Value = 44778
$ws.Range("E124").Value = 10
$ws.Range("F124").Value = 100112039
$ws.Range("G124").Value = "Ciboulette"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 240
$ws.Range("K124").Value = 4000
$ws.Range("L124").Value = 4000
$ws.Range("M124").Value = 4000
$ws.Range("N124").Value = "$/docena de atados"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 1333
$ws.Range("Q124").Value = 3
$ws.Range("R124").Value = "Hortaliza"
